# Simplify the "Instrument" sheet down to 5 generic instrument categories.
# (Commit: "Removed unnecessary 'Randomizer' file. Created an initial stat
#  generator for Members" - the Instrument list used by the randomizer is
#  trimmed from 21 specific instruments down to 5 broad categories.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Instrument")

# Remove the now-unused rows (6-21), keeping only rows 1-5.
$ws.Range("A6:A21").EntireRow.Delete()

# Replace the remaining 5 rows with the new, simplified category names.
$ws.Range("A1").Value = "Guitar"
$ws.Range("A2").Value = "Bass"
$ws.Range("A3").Value = "Drums"
$ws.Range("A4").Value = "Keyboard"
$ws.Range("A5").Value = "Vocals"
